$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 133.73334
$ws.Range("I42").Value = 114
$ws.Range("J42").Value = 173.2
$ws.Range("K42").Value = 342
$ws.Range("L42").Value = 519.5999999999999
$ws.Range("M42").Value = -112
$ws.Range("N42").Value = -979.5999999999999
$ws.Range("H62").Value = 4861.5454
$ws.Range("I62").Value = 2710.1428
$ws.Range("K62").Value = 2710.1428
$ws.Range("M62").Value = -2086.1428
$ws.Range("H65").Value = 4861.5454
$ws.Range("I65").Value = 2710.1428
$ws.Range("K65").Value = 13550.714
$ws.Range("M65").Value = -10430.714
$ws.Range("H82").Value = 525
$ws.Range("I82").Value = 525
$ws.Range("K82").Value = 1575
$ws.Range("M82").Value = -1169
$ws.Range("H85").Value = 525
$ws.Range("I85").Value = 525
$ws.Range("K85").Value = 1575
$ws.Range("M85").Value = -171
$ws.Range("H113").Value = 4885.4546
$ws.Range("J113").Value = 5224.1
$ws.Range("L113").Value = 5224.1
$ws.Range("N113").Value = -11732.1
$ws.Range("H132").Value = 1219.1428
$ws.Range("I132").Value = 929
$ws.Range("K132").Value = 2787
$ws.Range("M132").Value = -257
$ws.Range("H138").Value = 2774.1445
$ws.Range("J138").Value = 2951.0405
$ws.Range("L138").Value = 8853.121500000001
$ws.Range("N138").Value = -19133.1215
$ws.Range("H139").Value = 69999
$ws.Range("J139").Value = 69999
$ws.Range("L139").Value = 69999
$ws.Range("N139").Value = -80279
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4434.44
$ws.Range("I2").Value = 350.5
$ws.Range("K2").Value = 350.5
$ws.Range("M2").Value = -237.5
$ws.Range("H32").Value = 5104.222
$ws.Range("I32").Value = 3929.75
$ws.Range("K32").Value = 3929.75
$ws.Range("M32").Value = -3642.75
$ws.Range("H45").Value = 2742.7144
$ws.Range("I45").Value = 3041.5
$ws.Range("K45").Value = 3041.5
$ws.Range("M45").Value = -2664.5
$ws.Range("H116").Value = 4434.44
$ws.Range("I116").Value = 350.5
$ws.Range("K116").Value = 350.5
$ws.Range("M116").Value = 1943.5
$ws.Range("H122").Value = 5056.8887
$ws.Range("I122").Value = 6506
$ws.Range("J122").Value = 4642.857
$ws.Range("K122").Value = 19518
$ws.Range("L122").Value = 13928.571
$ws.Range("M122").Value = -17068
$ws.Range("N122").Value = -18828.571
$ws.Range("H132").Value = 2508.6072
$ws.Range("I132").Value = 1777.9615
$ws.Range("K132").Value = 5333.8845
$ws.Range("M132").Value = -2803.8845
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4434.44
$ws.Range("I3").Value = 350.5
$ws.Range("K3").Value = 350.5
$ws.Range("M3").Value = -236.5
$ws.Range("H64").Value = 1539.8
$ws.Range("I64").Value = 1539.8
$ws.Range("K64").Value = 1539.8
$ws.Range("M64").Value = -1314.8
$ws.Range("H67").Value = 1539.8
$ws.Range("I67").Value = 1539.8
$ws.Range("K67").Value = 1539.8
$ws.Range("M67").Value = -759.8
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 2212.125
$ws.Range("I22").Value = 589.8
$ws.Range("J22").Value = 4916
$ws.Range("K22").Value = 589.8
$ws.Range("L22").Value = 4916
$ws.Range("M22").Value = -239.8
$ws.Range("N22").Value = -5616
$ws.Range("H122").Value = 5581.353
$ws.Range("I122").Value = 2257.3333
$ws.Range("J122").Value = 13559
$ws.Range("K122").Value = 6771.999899999999
$ws.Range("L122").Value = 40677
$ws.Range("M122").Value = -4321.999899999999
$ws.Range("N122").Value = -45577
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H54").Value = 4166.6665
$ws.Range("J54").Value = 4250
$ws.Range("L54").Value = 12750
$ws.Range("N54").Value = -13868
$ws.Range("H119").Value = 7926.5713
$ws.Range("I119").Value = 2134.2856
$ws.Range("K119").Value = 6402.8568
$ws.Range("M119").Value = -1564.8568
$ws.Range("H136").Value = 41670276
$ws.Range("I136").Value = 55558704
$ws.Range("K136").Value = 166676112
$ws.Range("M136").Value = -166671012
$ws.Range("H140").Value = 1372.4348
$ws.Range("I140").Value = 1372.4348
$ws.Range("K140").Value = 4117.3044
$ws.Range("M140").Value = 1062.6956
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H34").Value = 23000
$ws.Range("J34").Value = 23000
$ws.Range("L34").Value = 23000
$ws.Range("N34").Value = -23536
$ws.Range("H76").Value = 23000
$ws.Range("J76").Value = 23000
$ws.Range("L76").Value = 23000
$ws.Range("N76").Value = -23630
$ws.Range("H79").Value = 23000
$ws.Range("J79").Value = 23000
$ws.Range("L79").Value = 23000
$ws.Range("N79").Value = -25184
$ws.Range("H102").Value = 2101.449
$ws.Range("I102").Value = 1683.2778
$ws.Range("J102").Value = 3259.4614
$ws.Range("K102").Value = 1683.2778
$ws.Range("L102").Value = 3259.4614
$ws.Range("M102").Value = -61.27780000000007
$ws.Range("N102").Value = -6503.4614
$ws.Range("H122").Value = 15049.23
$ws.Range("I122").Value = 23248
$ws.Range("J122").Value = 8021.7144
$ws.Range("K122").Value = 69744
$ws.Range("L122").Value = 24065.1432
$ws.Range("M122").Value = -67294
$ws.Range("N122").Value = -28965.1432
$ws.Range("H126").Value = 4571.684
$ws.Range("I126").Value = 3057.625
$ws.Range("J126").Value = 5672.8184
$ws.Range("K126").Value = 9172.875
$ws.Range("L126").Value = 17018.4552
$ws.Range("M126").Value = -6702.875
$ws.Range("N126").Value = -21958.4552
$ws.Range("H132").Value = 3970.3076
$ws.Range("I132").Value = 2985.8572
$ws.Range("J132").Value = 5118.8335
$ws.Range("K132").Value = 8957.571599999999
$ws.Range("L132").Value = 15356.5005
$ws.Range("M132").Value = -6427.571599999999
$ws.Range("N132").Value = -20416.5005
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 19845.818
$ws.Range("I40").Value = 28584
$ws.Range("J40").Value = 9360
$ws.Range("K40").Value = 28584
$ws.Range("L40").Value = 9360
$ws.Range("M40").Value = -28448
$ws.Range("N40").Value = -9632
$ws.Range("H46").Value = 2555.5642
$ws.Range("I46").Value = 1833.75
$ws.Range("J46").Value = 3057.6956
$ws.Range("K46").Value = 1833.75
$ws.Range("L46").Value = 3057.6956
$ws.Range("M46").Value = -1645.75
$ws.Range("N46").Value = -3433.6956
$ws.Range("H122").Value = 315316.7
$ws.Range("I122").Value = 1335066.6
$ws.Range("K122").Value = 4005199.8
$ws.Range("M122").Value = -4002749.8
$ws.Range("H136").Value = 4132.4727
$ws.Range("I136").Value = 3240.1843
$ws.Range("K136").Value = 9720.552899999999
$ws.Range("M136").Value = -7170.552899999999
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 49997.5
$ws.Range("J82").Value = 49996.668
$ws.Range("L82").Value = 49996.668
$ws.Range("N82").Value = -50762.668
$ws.Range("H85").Value = 49997.5
$ws.Range("J85").Value = 49996.668
$ws.Range("L85").Value = 49996.668
$ws.Range("N85").Value = -50762.668
$ws.Range("H86").Value = 49000
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("H89").Value = 49000
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("H113").Value = 436.25583
$ws.Range("I113").Value = 327.65518
$ws.Range("K113").Value = 982.9655399999999
$ws.Range("M113").Value = 1187.03446
$ws.Range("H132").Value = 7930.5
$ws.Range("I132").Value = 4405.6665
$ws.Range("K132").Value = 13216.9995
$ws.Range("M132").Value = -10686.9995
